$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get rotated/shuffled between rows (2..41)
$cols = @("D","L","M","N","O","P","Q","R","S","T")

# Mapping of new-row -> old-row that supplies its new D/L/M/N/O/P/Q/R/S/T values.
# (Derived by matching the post-edit row contents against the pre-edit rows.)
$rowMap = @{
    2  = 4
    3  = 5
    4  = 7
    5  = 28
    6  = 9
    7  = 32
    8  = 21
    9  = 33
    10 = 3
    11 = 30
    12 = 24
    13 = 15
    14 = 8
    15 = 23
    16 = 31
    17 = 40
    18 = 41
    19 = 27
    20 = 22
    21 = 26
    22 = 39
    23 = 19
    24 = 10
    25 = 34
    26 = 20
    27 = 12
    28 = 13
    29 = 14
    30 = 36
    31 = 25
    32 = 35
    33 = 6
    34 = 37
    35 = 2
    36 = 17
    37 = 18
    38 = 11
    39 = 16
    40 = 29
    41 = 38
}

# Snapshot all the original values first, since rows get cross-assigned.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Now write back the shuffled values.
for ($r = 2; $r -le 41; $r++) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
